$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3416.3333
$ws.Range("I43").Value = 3859.8
$ws.Range("K43").Value = 3859.8
$ws.Range("M43").Value = -3790.8
$ws.Range("H58").Value = 538.75
$ws.Range("I58").Value = 518.3333
$ws.Range("J58").Value = 600
$ws.Range("K58").Value = 1554.9999
$ws.Range("L58").Value = 1800
$ws.Range("M58").Value = -1404.9999
$ws.Range("N58").Value = -2100
$ws.Range("H70").Value = 4550.8887
$ws.Range("J70").Value = 5681.8887
$ws.Range("L70").Value = 17045.6661
$ws.Range("N70").Value = -17585.6661
$ws.Range("H73").Value = 4550.8887
$ws.Range("J73").Value = 5681.8887
$ws.Range("L73").Value = 17045.6661
$ws.Range("N73").Value = -18917.6661
$ws.Range("H100").Value = 2354.8667
$ws.Range("I100").Value = 2186.4614
$ws.Range("K100").Value = 2186.4614
$ws.Range("M100").Value = -1645.4614
$ws.Range("H111").Value = 1398.2
$ws.Range("I111").Value = 997.75
$ws.Range("K111").Value = 2993.25
$ws.Range("M111").Value = 73.75
$ws.Range("H138").Value = 1544
$ws.Range("I138").Value = 671.3570999999999
$ws.Range("J138").Value = 4598.25
$ws.Range("K138").Value = 2014.0713
$ws.Range("L138").Value = 13794.75
$ws.Range("M138").Value = 3125.9287
$ws.Range("N138").Value = -24074.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 3000
$ws.Range("K61").Value = 3000
$ws.Range("M61").Value = -2788
$ws.Range("H74").Value = 1425
$ws.Range("I74").Value = 1067.3334
$ws.Range("K74").Value = 1067.3334
$ws.Range("M74").Value = -193.3334
$ws.Range("H77").Value = 1425
$ws.Range("I77").Value = 1067.3334
$ws.Range("K77").Value = 5336.666999999999
$ws.Range("M77").Value = -968.6669999999995
$ws.Range("H102").Value = 1731.5333
$ws.Range("I102").Value = 1690.3077
$ws.Range("J102").Value = 1999.5
$ws.Range("K102").Value = 1690.3077
$ws.Range("L102").Value = 1999.5
$ws.Range("M102").Value = -68.30770000000007
$ws.Range("N102").Value = -5243.5
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 700
$ws.Range("I5").Value = 400
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 400
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -287
$ws.Range("N5").Value = -1226
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("N61").Value = 0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4445043
$ws.Range("J22").Value = 8000554
$ws.Range("L22").Value = 8000554
$ws.Range("N22").Value = -8001254
$ws.Range("H31").Value = 1831.1428
$ws.Range("I31").Value = 1869.8334
$ws.Range("K31").Value = 1869.8334
$ws.Range("M31").Value = -1574.8334
$ws.Range("H34").Value = 1831.1428
$ws.Range("I34").Value = 1869.8334
$ws.Range("K34").Value = 1869.8334
$ws.Range("M34").Value = -1667.8334
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1160.4
$ws.Range("I5").Value = 1046.7142
$ws.Range("K5").Value = 3140.1426
$ws.Range("M5").Value = -3028.1426
$ws.Range("H135").Value = 1160.4
$ws.Range("I135").Value = 1046.7142
$ws.Range("K135").Value = 9420.427799999999
$ws.Range("M135").Value = -6885.427799999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 116.28571
$ws.Range("I2").Value = 134
$ws.Range("K2").Value = 134
$ws.Range("M2").Value = -21
$ws.Range("H3").Value = 506000
$ws.Range("I3").Value = 506000
$ws.Range("K3").Value = 506000
$ws.Range("M3").Value = -505884
$ws.Range("H59").Value = 12000
$ws.Range("J59").Value = 12000
$ws.Range("L59").Value = 12000
$ws.Range("N59").Value = -13166
$ws.Range("H70").Value = 9866
$ws.Range("I70").Value = 10225.728
$ws.Range("K70").Value = 10225.728
$ws.Range("M70").Value = -9955.727999999999
$ws.Range("H73").Value = 9866
$ws.Range("I73").Value = 10225.728
$ws.Range("K73").Value = 10225.728
$ws.Range("M73").Value = -9289.727999999999
$ws.Range("H80").Value = 4832.6665
$ws.Range("J80").Value = 2500
$ws.Range("L80").Value = 2500
$ws.Range("N80").Value = -4496
$ws.Range("H83").Value = 4832.6665
$ws.Range("J83").Value = 2500
$ws.Range("L83").Value = 12500
$ws.Range("N83").Value = -22484
$ws.Range("H122").Value = 2579.35
$ws.Range("I122").Value = 2112.375
$ws.Range("K122").Value = 6337.125
$ws.Range("M122").Value = -3887.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7155
$ws.Range("I7").Value = 6996.5454
$ws.Range("J7").Value = 8898
$ws.Range("K7").Value = 6996.5454
$ws.Range("L7").Value = 6996.5454
$ws.Range("M7").Value = -6884.5454
$ws.Range("N7").Value = -9122
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 5000
$ws.Range("K25").Value = 5000
$ws.Range("M25").Value = -4770
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H122").Value = 4043.1052
$ws.Range("I122").Value = 3618.5
$ws.Range("K122").Value = 10855.5
$ws.Range("M122").Value = -8405.5
$ws.Range("H126").Value = 7155
$ws.Range("I126").Value = 6996.5454
$ws.Range("J126").Value = 8898
$ws.Range("K126").Value = 20989.6362
$ws.Range("L126").Value = 26694
$ws.Range("M126").Value = -18519.6362
$ws.Range("N126").Value = -31634
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12696.385
$ws.Range("I62").Value = 10119.8
$ws.Range("J62").Value = 14306.75
$ws.Range("K62").Value = 10119.8
$ws.Range("L62").Value = 14306.75
$ws.Range("M62").Value = -9495.799999999999
$ws.Range("N62").Value = -15554.75
$ws.Range("H65").Value = 12696.385
$ws.Range("I65").Value = 10119.8
$ws.Range("J65").Value = 14306.75
$ws.Range("K65").Value = 50599
$ws.Range("L65").Value = 71533.75
$ws.Range("M65").Value = -47479
$ws.Range("N65").Value = -77773.75
